$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data, sorted by value descending, with Swedish and Uzbek removed.
$data = @(
    @("Chinese", 21.0249771737996),
    @("English", 19.9005116352777),
    @("Spanish", 5.871907646539979),
    @("Arabic", 4.264745979855292),
    @("German", 3.688852290967259),
    @("Malay-Indonesian", 3.511511598918445),
    @("Japanese", 3.275974609939562),
    @("Russian", 2.686691419322029),
    @("Portuguese", 2.628992721870851),
    @("French", 2.340479079820591),
    @("Turkish", 2.05345554890109),
    @("Italian", 1.65610064261636),
    @("Korean", 1.614256154304532),
    @("Dutch", 1.092162329866998),
    @("Bengali", 1.020985449591072),
    @("Vietnamese", 0.9900823273187158),
    @("Polish", 0.9817816872266322),
    @("Urdu", 0.9744646011000273),
    @("Persian", 0.9228965442132723),
    @("Thai", 0.9184902582786737)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove the now-unused trailing rows (previously rows 22 and 23).
$ws.Range("A22:B23").Delete()
